# Journal de travail - François
# Replace the placeholder "bla bla" activity text with the real activity
# descriptions, and correct the hours logged for each entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: first activity + hours
$ws.Range("B5").Value = "Choix de la proposition du projet "
$ws.Range("C5").Value = 0.25

# Row 6: second activity + hours
$ws.Range("B6").Value = "Lecture de la proposition"
$ws.Range("C6").Value = 0.25
